$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add column G across all 16 rows (same "default" cell style as column A/B/F) ---
for ($r = 1; $r -le 16; $r++) {
    $ws.Cells.Item($r, 2).Copy()
    $ws.Cells.Item($r, 7).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# --- Row 13 (pin PA4) content updates ---
# "说明" (description) changes
$ws.Range("D13").Value = "FSK控制的慢1/4相位"
# New "用途" (usage) value placed into the newly-added column G
$ws.Range("G13").Value = "PWM输出"

# --- Normalize the merged-cell border styling on F7/F16 to match F6/F15 ---
$ws.Range("F6").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F15").Copy()
$ws.Range("F16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- View state: zoom level and active selection ---
$ws.Application.ActiveWindow.Zoom = 190
$ws.Range("D18").Select() | Out-Null
